$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.611809333333333
$ws.Range("H2").Value = 22.835428
$ws.Range("I2").Value = 0.1134603829630287
$ws.Range("J2").Value = 0.1134603829630287
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.110028
$ws.Range("N2").Value = 0.330084
$ws.Range("Q2").Value = 0.837512157328
$ws.Range("R2").Value = 7.537609415952
$ws.Range("S2").Value = 0.1134603829630287
$ws.Range("T2").Value = 0.1134603829630287

$ws.Range("I3").Value = 0.1569061007197586
$ws.Range("J3").Value = 0.1569061007197586
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.110028
$ws.Range("N3").Value = 0.330084
$ws.Range("Q3").Value = 1.158208384988
$ws.Range("R3").Value = 10.423875464892
$ws.Range("S3").Value = 0.1569061007197586
$ws.Range("T3").Value = 0.1569061007197586

$ws.Range("G4").Value = 19.26381566666667
$ws.Range("H4").Value = 57.791447
$ws.Range("I4").Value = 0.2871432805466829
$ws.Range("J4").Value = 0.2871432805466829
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110028
$ws.Range("N4").Value = 0.330084
$ws.Range("Q4").Value = 2.119559110172
$ws.Range("R4").Value = 19.076031991548
$ws.Range("S4").Value = 0.2871432805466829
$ws.Range("T4").Value = 0.2871432805466829

$ws.Range("G5").Value = 3.278219666666667
$ws.Range("H5").Value = 9.834659
$ws.Range("I5").Value = 0.04886460531638807
$ws.Range("J5").Value = 0.04886460531638808
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.110028
$ws.Range("N5").Value = 0.330084
$ws.Range("Q5").Value = 0.360695953484
$ws.Range("R5").Value = 3.246263581356
$ws.Range("S5").Value = 0.04886460531638807
$ws.Range("T5").Value = 0.04886460531638808

$ws.Range("G6").Value = 5.160004333333333
$ws.Range("H6").Value = 15.480013
$ws.Range("I6").Value = 0.07691417928547969
$ws.Range("J6").Value = 0.07691417928547971
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.110028
$ws.Range("N6").Value = 0.330084
$ws.Range("Q6").Value = 0.567744956788
$ws.Range("R6").Value = 5.109704611092
$ws.Range("S6").Value = 0.07691417928547969
$ws.Range("T6").Value = 0.07691417928547971

$ws.Range("G7").Value = 21.24747966666667
$ws.Range("H7").Value = 63.742439
$ws.Range("I7").Value = 0.316711451168662
$ws.Range("J7").Value = 0.3167114511686621
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.110028
$ws.Range("N7").Value = 0.330084
$ws.Range("Q7").Value = 2.337817692764
$ws.Range("R7").Value = 21.040359234876
$ws.Range("S7").Value = 0.316711451168662
$ws.Range("T7").Value = 0.3167114511686621
